$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").ClearContents()
$ws.Range("C4").Value = "Jadot"
$ws.Range("D4").Value = "Jadot"
$ws.Range("E4").Value = "Jadot"
$ws.Range("F4").Value = "Jadot"
$ws.Range("G4").Value = "Jadot"
$ws.Range("H4").Value = "Macron"
$ws.Range("I4").Value = "Jadot"
$ws.Range("J4").Value = "Macron"
$ws.Range("K4").Value = "Macron"
$ws.Range("L4").Value = "Macron"
$ws.Range("M4").Value = "Macron"
$ws.Range("N4").Value = "Macron"
$ws.Range("O4").Value = "Philipot"
$ws.Range("Q4").Value = "Philipot"
$ws.Range("R4").Value = "Philipot"
$ws.Range("S4").Value = "Philipot"
$ws.Range("U4").Value = "Philipot"
$ws.Range("W4").Value = "Philipot"
$ws.Range("X4").Value = "Philipot"
$ws.Range("Y4").Value = "Philipot"
$ws.Range("AA4").Value = "Philipot"
$ws.Range("AB4").Value = "Montebourg"
$ws.Range("AG4").Value = "Philipot"
$ws.Range("AH4").Value = "Montebourg"
$ws.Range("AM4").Value = "Philipot"
$ws.Range("AN4").Value = "Montebourg"
$ws.Range("B5").Value = "Jadot"
$ws.Range("C5").Value = "Jadot"
$ws.Range("D5").Value = "Jadot"
$ws.Range("E5").Value = "Jadot"
$ws.Range("F5").Value = "Jadot"
$ws.Range("G5").Value = "Jadot"
$ws.Range("H5").Value = "Macron"
$ws.Range("I5").Value = "Jadot"
$ws.Range("J5").Value = "Jadot"
$ws.Range("K5").Value = "Jadot"
$ws.Range("L5").Value = "Jadot"
$ws.Range("M5").Value = "Jadot"
$ws.Range("N5").Value = "Macron"
$ws.Range("O5").Value = "Jadot"
$ws.Range("P5").Value = "Jadot"
$ws.Range("Q5").Value = "Jadot"
$ws.Range("R5").Value = "Jadot"
$ws.Range("S5").Value = "Jadot"
$ws.Range("V5").Value = "Philipot"
$ws.Range("Y5").Value = "Mélenchon"
$ws.Range("B6").Value = "Jadot"
$ws.Range("C6").Value = "Jadot"
$ws.Range("D6").Value = "Jadot"
$ws.Range("E6").Value = "Jadot"
$ws.Range("F6").Value = "Jadot"
$ws.Range("G6").Value = "Jadot"
$ws.Range("H6").Value = "Macron"
$ws.Range("I6").Value = "Jadot"
$ws.Range("J6").Value = "Macron"
$ws.Range("K6").Value = "Macron"
$ws.Range("L6").Value = "Macron"
$ws.Range("M6").Value = "Jadot"
$ws.Range("O6").Value = "Jadot"
$ws.Range("Q6").Value = "Macron"
$ws.Range("R6").Value = "Jadot"
$ws.Range("S6").Value = "Jadot"
$ws.Range("T6").Value = "Philipot"
$ws.Range("U6").Value = "Jadot"
$ws.Range("V6").Value = "Philipot"
$ws.Range("W6").Value = "Philipot"
$ws.Range("X6").Value = "Macron"
$ws.Range("Y6").Value = "Mélenchon"
$ws.Range("Z6").Value = "Philipot"
$ws.Range("AA6").Value = "Jadot"
$ws.Range("AF6").Value = "Philipot"
$ws.Range("AG6").Value = "Jadot"
$ws.Range("AL6").Value = "Philipot"
$ws.Range("AM6").Value = "Philipot"
$ws.Range("AN6").Value = "Philipot"
$ws.Range("B7").Value = "Jadot"
$ws.Range("C7").Value = "Jadot"
$ws.Range("D7").Value = "Jadot"
$ws.Range("E7").Value = "Jadot"
$ws.Range("F7").Value = "Jadot"
$ws.Range("G7").Value = "Jadot"
$ws.Range("H7").Value = "Macron"
$ws.Range("I7").Value = "Macron"
$ws.Range("J7").Value = "Macron"
$ws.Range("K7").Value = "Macron"
$ws.Range("L7").Value = "Macron"
$ws.Range("M7").Value = "Macron"
$ws.Range("N7").Value = "Macron"
$ws.Range("O7").Value = "Philipot"
$ws.Range("P7").Value = "Macron"
$ws.Range("Q7").Value = "Philipot"
$ws.Range("R7").Value = "Philipot"
$ws.Range("S7").Value = "Philipot"
$ws.Range("U7").Value = "Philipot"
$ws.Range("V7").Value = "Philipot"
$ws.Range("W7").Value = "Philipot"
$ws.Range("X7").Value = "Philipot"
$ws.Range("Y7").Value = "Philipot"
$ws.Range("Z7").Value = "Philipot"
$ws.Range("AA7").Value = "Philipot"
$ws.Range("AF7").Value = "Macron"
$ws.Range("AG7").Value = "Macron"
$ws.Range("AH7").Value = "Macron"
$ws.Range("AL7").Value = "Philipot"
$ws.Range("AN7").Value = "Philipot"
$ws.Range("B8").Value = "Jadot"
$ws.Range("C8").Value = "Jadot"
$ws.Range("D8").Value = "Jadot"
$ws.Range("E8").Value = "Jadot"
$ws.Range("F8").Value = "Jadot"
$ws.Range("G8").Value = "Jadot"
$ws.Range("H8").Value = "Macron"
$ws.Range("I8").Value = "Macron"
$ws.Range("J8").Value = "Macron"
$ws.Range("K8").Value = "Macron"
$ws.Range("L8").Value = "Macron"
$ws.Range("M8").Value = "Macron"
$ws.Range("N8").Value = "Macron"
$ws.Range("O8").Value = "Philipot"
$ws.Range("P8").Value = "Philipot"
$ws.Range("Q8").Value = "Philipot"
$ws.Range("R8").Value = "Philipot"
$ws.Range("S8").Value = "Philipot"
$ws.Range("T8").Value = "Macron"
$ws.Range("U8").Value = "Philipot"
$ws.Range("V8").Value = "Philipot"
$ws.Range("W8").Value = "Philipot"
$ws.Range("X8").Value = "Philipot"
$ws.Range("Y8").Value = "Philipot"
$ws.Range("Z8").Value = "Philipot"
$ws.Range("AA8").Value = "Philipot"
$ws.Range("AF8").Value = "Macron"
$ws.Range("AG8").Value = "Philipot"
$ws.Range("AL8").Value = "Macron"
$ws.Range("AM8").Value = "Philipot"
$ws.Range("B9").Value = "Jadot"
$ws.Range("C9").Value = "Jadot"
$ws.Range("D9").Value = "Jadot"
$ws.Range("E9").Value = "Jadot"
$ws.Range("F9").Value = "Jadot"
$ws.Range("G9").Value = "Jadot"
$ws.Range("H9").Value = "Macron"
$ws.Range("I9").Value = "Jadot"
$ws.Range("J9").Value = "Jadot"
$ws.Range("K9").Value = "Macron"
$ws.Range("L9").Value = "Jadot"
$ws.Range("M9").Value = "Jadot"
$ws.Range("N9").Value = "Macron"
$ws.Range("O9").Value = "Jadot"
$ws.Range("Q9").Value = "Macron"
$ws.Range("R9").Value = "Jadot"
$ws.Range("S9").Value = "Jadot"
$ws.Range("U9").Value = "Jadot"
$ws.Range("V9").Value = "Philipot"
$ws.Range("W9").Value = "Macron"
$ws.Range("X9").Value = "Jadot"
$ws.Range("Y9").Value = "Jadot"
$ws.Range("AA9").Value = "Jadot"
$ws.Range("AG9").Value = "Jadot"
$ws.Range("AH9").Value = "Jadot"
$ws.Range("AK9").Value = "Montebourg"
$ws.Range("AM9").Value = "Jadot"
$ws.Range("B10").Value = "Zemmour"
$ws.Range("C10").Value = "Jadot"
$ws.Range("D10").Value = "Jadot"
$ws.Range("E10").Value = "Jadot"
$ws.Range("F10").Value = "Jadot"
$ws.Range("G10").Value = "Jadot"
$ws.Range("H10").Value = "Macron"
$ws.Range("I10").Value = "Macron"
$ws.Range("J10").Value = "Macron"
$ws.Range("K10").Value = "Macron"
$ws.Range("L10").Value = "Macron"
$ws.Range("M10").Value = "Macron"
$ws.Range("N10").Value = "Macron"
$ws.Range("Q10").Value = "Philipot"
$ws.Range("R10").Value = "Philipot"
$ws.Range("S10").Value = "Philipot"
$ws.Range("T10").Value = "Macron"
$ws.Range("W10").Value = "Philipot"
$ws.Range("X10").Value = "Philipot"
$ws.Range("Y10").Value = "Philipot"
$ws.Range("Z10").Value = "Macron"
$ws.Range("AB10").Value = "Montebourg"
$ws.Range("AF10").Value = "Macron"
$ws.Range("AH10").Value = "Montebourg"
$ws.Range("AL10").Value = "Macron"
$ws.Range("AN10").Value = "Montebourg"
$ws.Range("B11").Value = "Jadot"
$ws.Range("C11").Value = "Jadot"
$ws.Range("D11").Value = "Jadot"
$ws.Range("E11").Value = "Jadot"
$ws.Range("F11").Value = "Jadot"
$ws.Range("G11").Value = "Jadot"
$ws.Range("H11").Value = "Macron"
$ws.Range("I11").Value = "Macron"
$ws.Range("J11").Value = "Macron"
$ws.Range("K11").Value = "Macron"
$ws.Range("L11").Value = "Macron"
$ws.Range("M11").Value = "Macron"
$ws.Range("N11").Value = "Macron"
$ws.Range("O11").Value = "Philipot"
$ws.Range("P11").Value = "Philipot"
$ws.Range("Q11").Value = "Philipot"
$ws.Range("R11").Value = "Philipot"
$ws.Range("S11").Value = "Philipot"
$ws.Range("U11").Value = "Philipot"
$ws.Range("V11").Value = "Philipot"
$ws.Range("W11").Value = "Philipot"
$ws.Range("X11").Value = "Philipot"
$ws.Range("Y11").Value = "Philipot"
$ws.Range("B12").Value = "Jadot"
$ws.Range("C12").Value = "Jadot"
$ws.Range("D12").Value = "Jadot"
$ws.Range("E12").Value = "Jadot"
$ws.Range("F12").Value = "Jadot"
$ws.Range("G12").Value = "Jadot"
$ws.Range("H12").Value = "Macron"
$ws.Range("I12").Value = "Macron"
$ws.Range("J12").Value = "Macron"
$ws.Range("K12").Value = "Macron"
$ws.Range("L12").Value = "Macron"
$ws.Range("M12").Value = "Macron"
$ws.Range("N12").Value = "Macron"
$ws.Range("O12").Value = "Philipot"
$ws.Range("P12").Value = "Philipot"
$ws.Range("Q12").Value = "Philipot"
$ws.Range("R12").Value = "Philipot"
$ws.Range("S12").Value = "Philipot"
$ws.Range("T12").Value = "Macron"
$ws.Range("U12").Value = "Philipot"
$ws.Range("V12").Value = "Philipot"
$ws.Range("W12").Value = "Philipot"
$ws.Range("X12").Value = "Philipot"
$ws.Range("Y12").Value = "Philipot"
$ws.Range("Z12").Value = "Philipot"
$ws.Range("AA12").Value = "Philipot"
$ws.Range("AF12").Value = "Macron"
$ws.Range("AG12").Value = "Philipot"
$ws.Range("AL12").Value = "Macron"
$ws.Range("AM12").Value = "Philipot"
